# Implements a small "writeline" helper that writes a row of values to a
# worksheet (one cell per column letter supplied), always storing the
# values as text - mirroring how a CSV/record writer would dump a line of
# fields into a sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Write-ExcelLine($Row, $Cols, $Values) {
    for ($i = 0; $i -lt $Cols.Length; $i++) {
        $addr = "$($Cols[$i])$Row"
        $cell = $ws.Range($addr)

        # Force text storage (so "123"/"1"/"0" aren't coerced to numbers),
        # then drop the number-format override again so the cell is left
        # with no extra styling - just the text value.
        $cell.NumberFormat = "@"
        $cell.Value = [string]$Values[$i]
        $cell.ClearFormats()
    }
}

# Remove the previous sample data that used to live in the sheet.
$ws.Range("A1").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("A4").ClearContents()

# Write the header/record line.
$columns = @("A", "B", "C", "D", "F", "H", "J", "L", "N", "P")
$values  = @("Dr. John Doe", "123", "2022-06-04", "1", "0", "0", "0", "0", "1", "0")

Write-ExcelLine 1 $columns $values
